$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1850.4286
$ws.Range("I19").Value = 1406.0714
$ws.Range("K19").Value = 1406.0714
$ws.Range("M19").Value = -1231.0714
$ws.Range("H32").Value = 6817.231
$ws.Range("I32").Value = 7311
$ws.Range("J32").Value = 6241.1665
$ws.Range("K32").Value = 7311
$ws.Range("L32").Value = 6241.1665
$ws.Range("M32").Value = -6985
$ws.Range("N32").Value = -6893.1665
$ws.Range("H87").Value = 99160.336
$ws.Range("J87").Value = 99160.336
$ws.Range("L87").Value = 99160.336
$ws.Range("N87").Value = -101656.336
$ws.Range("H90").Value = 99160.336
$ws.Range("J90").Value = 99160.336
$ws.Range("L90").Value = 297481.008
$ws.Range("N90").Value = -309961.008
$ws.Range("H132").Value = 2276.3923
$ws.Range("I132").Value = 2268.25
$ws.Range("K132").Value = 6804.75
$ws.Range("M132").Value = -4274.75
$ws.Range("H135").Value = 5578.64
$ws.Range("I135").Value = 1691.85
$ws.Range("K135").Value = 15226.65
$ws.Range("M135").Value = -12691.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 24071.666
$ws.Range("I45").Value = 26143.125
$ws.Range("K45").Value = 26143.125
$ws.Range("M45").Value = -25766.125
$ws.Range("H61").Value = 4239.4736
$ws.Range("I61").Value = 3837.8572
$ws.Range("J61").Value = 8925
$ws.Range("K61").Value = 3837.8572
$ws.Range("L61").Value = 8925
$ws.Range("M61").Value = -3625.8572
$ws.Range("N61").Value = -9349
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 109670
$ws.Range("I81").Value = 59999.5
$ws.Range("K81").Value = 59999.5
$ws.Range("M81").Value = -59001.5
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 109670
$ws.Range("I84").Value = 59999.5
$ws.Range("K84").Value = 179998.5
$ws.Range("M84").Value = -175006.5
$ws.Range("H102").Value = 4209.75
$ws.Range("I102").Value = 4175.391
$ws.Range("K102").Value = 4175.391
$ws.Range("M102").Value = -2553.391
$ws.Range("H110").Value = 1278.625
$ws.Range("I110").Value = 1381.9231
$ws.Range("J110").Value = 831
$ws.Range("K110").Value = 1381.9231
$ws.Range("L110").Value = 831
$ws.Range("M110").Value = 663.0769
$ws.Range("N110").Value = -4921
$ws.Range("H122").Value = 2624.6667
$ws.Range("I122").Value = 2197.2856
$ws.Range("K122").Value = 6591.8568
$ws.Range("M122").Value = -4141.8568
$ws.Range("H136").Value = 4239.4736
$ws.Range("I136").Value = 3837.8572
$ws.Range("J136").Value = 8925
$ws.Range("K136").Value = 11513.5716
$ws.Range("L136").Value = 26775
$ws.Range("M136").Value = -8963.571599999999
$ws.Range("N136").Value = -31875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 233.33333
$ws.Range("I22").Value = 233.33333
$ws.Range("K22").Value = 233.33333
$ws.Range("M22").Value = -60.33332999999999
$ws.Range("H86").Value = 3661.4
$ws.Range("I86").Value = 3514.5
$ws.Range("J86").Value = 4249
$ws.Range("K86").Value = 3514.5
$ws.Range("L86").Value = 4249
$ws.Range("M86").Value = -2391.5
$ws.Range("N86").Value = -6495
$ws.Range("H89").Value = 3661.4
$ws.Range("I89").Value = 3514.5
$ws.Range("J89").Value = 4249
$ws.Range("K89").Value = 17572.5
$ws.Range("L89").Value = 21245
$ws.Range("M89").Value = -11956.5
$ws.Range("N89").Value = -32477
$ws.Range("H99").Value = 3762.8
$ws.Range("I99").Value = 3638.7144
$ws.Range("K99").Value = 3638.7144
$ws.Range("M99").Value = -2140.7144
$ws.Range("H105").Value = 4144.2915
$ws.Range("I105").Value = 4165.1904
$ws.Range("K105").Value = 4165.1904
$ws.Range("M105").Value = -2418.1904
$ws.Range("H134").Value = 1817.6666
$ws.Range("I134").Value = 1817.6666
$ws.Range("K134").Value = 5452.9998
$ws.Range("M134").Value = -2917.9998
$ws.Range("H135").Value = 73332
$ws.Range("J135").Value = 73332
$ws.Range("L135").Value = 73332
$ws.Range("N135").Value = -83472
$ws.Range("H137").Value = 79760
$ws.Range("J137").Value = 79760
$ws.Range("L137").Value = 79760
$ws.Range("N137").Value = -89960
$ws.Range("H138").Value = 89602.28999999999
$ws.Range("J138").Value = 89602.28999999999
$ws.Range("L138").Value = 89602.28999999999
$ws.Range("N138").Value = -99882.28999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4565.615
$ws.Range("I31").Value = 3406.6667
$ws.Range("K31").Value = 3406.6667
$ws.Range("M31").Value = -3111.6667
$ws.Range("H34").Value = 4565.615
$ws.Range("I34").Value = 3406.6667
$ws.Range("K34").Value = 3406.6667
$ws.Range("M34").Value = -3204.6667
$ws.Range("H58").Value = 3950.9143
$ws.Range("I58").Value = 3894
$ws.Range("K58").Value = 3894
$ws.Range("M58").Value = -3691
$ws.Range("H105").Value = 1135.6364
$ws.Range("I105").Value = 1135.6364
$ws.Range("K105").Value = 1135.6364
$ws.Range("M105").Value = 611.3635999999999
$ws.Range("H135").Value = 99997
$ws.Range("J135").Value = 99997
$ws.Range("L135").Value = 99997
$ws.Range("N135").Value = -110137
$ws.Range("H136").Value = 3950.9143
$ws.Range("I136").Value = 3894
$ws.Range("K136").Value = 11682
$ws.Range("M136").Value = -9132
$ws.Range("H138").Value = 88599.39999999999
$ws.Range("J138").Value = 88599.39999999999
$ws.Range("L138").Value = 88599.39999999999
$ws.Range("N138").Value = -98879.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 125045624
$ws.Range("J37").Value = 125045624
$ws.Range("L37").Value = 375136872
$ws.Range("N37").Value = -375137096
$ws.Range("H68").Value = 1364.875
$ws.Range("I68").Value = 716.3333
$ws.Range("J68").Value = 1754
$ws.Range("K68").Value = 2148.9999
$ws.Range("L68").Value = 5262
$ws.Range("M68").Value = -1337.9999
$ws.Range("N68").Value = -6884
$ws.Range("H71").Value = 1364.875
$ws.Range("I71").Value = 716.3333
$ws.Range("J71").Value = 1754
$ws.Range("K71").Value = 6446.9997
$ws.Range("L71").Value = 15786
$ws.Range("M71").Value = -2390.9997
$ws.Range("N71").Value = -23898
$ws.Range("H137").Value = 3495.6553
$ws.Range("I137").Value = 3539.1538
$ws.Range("J137").Value = 3460.3125
$ws.Range("K137").Value = 10617.4614
$ws.Range("L137").Value = 10380.9375
$ws.Range("M137").Value = -5517.4614
$ws.Range("N137").Value = -20580.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5900.7
$ws.Range("I80").Value = 3002.5
$ws.Range("K80").Value = 3002.5
$ws.Range("M80").Value = -2004.5
$ws.Range("H83").Value = 5900.7
$ws.Range("I83").Value = 3002.5
$ws.Range("K83").Value = 15012.5
$ws.Range("M83").Value = -10020.5
$ws.Range("H133").Value = 79365
$ws.Range("J133").Value = 79365
$ws.Range("L133").Value = 79365
$ws.Range("N133").Value = -89485
$ws.Range("H135").Value = 84130
$ws.Range("J135").Value = 84130
$ws.Range("L135").Value = 84130
$ws.Range("N135").Value = -94270
$ws.Range("H138").Value = 69993
$ws.Range("J138").Value = 69993
$ws.Range("L138").Value = 69993
$ws.Range("N138").Value = -80273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5652.231
$ws.Range("J46").Value = 3118.8
$ws.Range("L46").Value = 3118.8
$ws.Range("N46").Value = -3494.8
$ws.Range("H82").Value = 983.8461
$ws.Range("I82").Value = 1013.25
$ws.Range("J82").Value = 936.8
$ws.Range("K82").Value = 1013.25
$ws.Range("L82").Value = 936.8
$ws.Range("M82").Value = -652.25
$ws.Range("N82").Value = -1658.8
$ws.Range("H85").Value = 983.8461
$ws.Range("I85").Value = 1013.25
$ws.Range("J85").Value = 936.8
$ws.Range("K85").Value = 1013.25
$ws.Range("L85").Value = 936.8
$ws.Range("M85").Value = 234.75
$ws.Range("N85").Value = -3432.8
$ws.Range("H125").Value = 73158.2
$ws.Range("J125").Value = 73158.2
$ws.Range("L125").Value = 73158.2
$ws.Range("N125").Value = -82998.2
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2794.4
$ws.Range("J107").Value = 3889.5
$ws.Range("L107").Value = 11668.5
$ws.Range("N107").Value = -15508.5
$ws.Range("H125").Value = 65712
$ws.Range("J125").Value = 65712
$ws.Range("L125").Value = 65712
$ws.Range("N125").Value = -75552
$ws.Range("H132").Value = 1887.7142
$ws.Range("I132").Value = 1393.4546
$ws.Range("K132").Value = 4180.3638
$ws.Range("M132").Value = -1650.3638
$ws.Range("H136").Value = 2964.2144
$ws.Range("I136").Value = 3320.8333
$ws.Range("K136").Value = 9962.499899999999
$ws.Range("M136").Value = -7412.499899999999
